$d = $word.ActiveDocument

# --- Simple "remove trailing year/date parenthetical" edits on bold heading runs ---
# These each occur twice in the document (once as a plain descriptive run and
# once as a bold heading run); only the bold heading run should change.

$d.Paragraphs(28).Range.Find.Execute(
    "Medietekniks handledningsprocesser 15hp (2006)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Medietekniks handledningsprocesser 15hp", 2)

$d.Paragraphs(30).Range.Find.Execute(
    "Handledning i högre utbildning (2014)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Handledning i högre utbildning", 2)

$d.Paragraphs(32).Range.Find.Execute(
    "Klart ledarskap (2021)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Klart ledarskap", 2)

$d.Paragraphs(35).Range.Find.Execute(
    "Universitetsajunkt @ BTH (JUN 2004 - pågående)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Universitetsajunkt @ BTH", 2)

$d.Paragraphs(60).Range.Find.Execute(
    "Back-end utvecklare @ ManagerZone (2007)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Back-end utvecklare @ ManagerZone", 2)

$d.Paragraphs(62).Range.Find.Execute(
    "Musiker @ Playtones (2010 - 2013)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Musiker @ Playtones", 2)

$d.Paragraphs(65).Range.Find.Execute(
    "Spider SMS (2005-2006)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Spider SMS", 2)

$d.Paragraphs(67).Range.Find.Execute(
    "Kombinerad mobilitet - Växla upp! (2018-2020)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Kombinerad mobilitet - Växla upp!", 2)

# --- Split long description runs into extra bulleted lines (manual line breaks) ---

# Musiker @ Playtones description paragraph (para 61)
$d.Paragraphs(61).Range.Find.Execute(
    "- Under tre år var jag tjänstledig från BTH och försöjde jag mig som musiker i gruppen The Playtones. Efter medverkan och vinst i TV-programmet Dansbandskampen vintern 2009, blev en hobby ett leverbröd.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "- Under tre år var jag tjänstledig från BTH och försöjde jag mig som musiker i gruppen The Playtones.^l- Efter medverkan och vinst i TV-programmet Dansbandskampen vintern 2009, blev en hobby ett leverbröd.",
    2)

$d.Paragraphs(61).Range.Find.Execute(
    "Intressant att få uppleva musikindustrin från insidan. Uppträden i TV-produktioner såsom Allsång på skansen och Melodifestivalen.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "- Intressant att få uppleva musikindustrin från insidan. Uppträden i TV-produktioner såsom Allsång på skansen och Melodifestivalen.",
    2)

# Spider SMS description paragraph (para 64)
$d.Paragraphs(64).Range.Find.Execute(
    "- Samarbete mellan BTH, Makerere University (Kampala, Uganda) och De La Salle University (Manila, Filippinerna). Delvis finansierat av medverkande partners samt SIDA. Projektets mål var att utveckla en applikation för bönder på landsbyggden i Uganda, för att via SMS få information om dagspriser på mejeriprodukter. Min roll i teamet var att ta fram en server-applikation för att lagra och handha priser på mejeriprodukter (API). Tekniker var html/css/javascript, php och MySql. Från BTH var professor Lena Trojer projektledare.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "- Samarbete mellan BTH, Makerere University (Kampala, Uganda) och De La Salle University (Manila, Filippinerna).^l- Delvis finansierat av medverkande partners samt SIDA. Projektets mål var att utveckla en applikation för bönder på landsbyggden i Uganda, för att via SMS få information om dagspriser på mejeriprodukter.^l- Min roll i teamet var att ta fram en server-applikation för att lagra och handha priser på mejeriprodukter (API). Tekniker var html/css/javascript, php och MySql. Från BTH var professor Lena Trojer projektledare.",
    2)

# Kombinerad mobilitet - Växla upp! description paragraph (para 66)
$d.Paragraphs(66).Range.Find.Execute(
    "- Projekt med det övergripande målet att minska koldioxidutsläppen från persontransporter. Finansierades av Europeiska regionala utvecklingsfonden (ERUF) där Energikontor sydost var projektägare.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "- Projekt med det övergripande målet att minska koldioxidutsläppen från persontransporter.^l- Finansierades av Europeiska regionala utvecklingsfonden (ERUF) där Energikontor sydost var projektägare.",
    2)

$d.Paragraphs(66).Range.Find.Execute(
    "- Min medverkan var i ett team från BTH och Netport som utvecklade appen Växla upp! (iPhone, Android) för att registrera cykelturer. Jag arbetade med att utveckla servern mottog, lagrade och analyserade registrerad data under användarnas cykelturer. Teknikerna var node/javascript, MongoDb, git och testramverket Jest.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "- Min medverkan var i ett team från BTH och Netport som utvecklade appen Växla upp! (iPhone, Android) för att registrera cykelturer.^l- Jag arbetade med att utveckla servern mottog, lagrade och analyserade registrerad data under användarnas cykelturer. Teknikerna var node/javascript, MongoDb, git och testramverket Jest.",
    2)

Write-Output "All edits applied"
